# Apply the Aug 31 2023 cryptos list refresh (prices + 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "26.581.59" or "1.009" that must stay text
# (Excel would otherwise silently coerce single-dot values like "219.58" into
# numbers and round-trip them with float noise). Force text format for the
# whole price column before writing, then restore the original (default) style
# so no stray formatting is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.549.91"
$ws.Range("E2").Value = "  -2.57%  "

# Row 3
$ws.Range("D3").Value = "1.668.92"
$ws.Range("E3").Value = "  -2.12%  "

# Row 4
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.61%  "

# Row 5
$ws.Range("D5").Value = "219.58"
$ws.Range("E5").Value = "  -1.58%  "

# Row 6
$ws.Range("D6").Value = "0.5127"
$ws.Range("E6").Value = "  -3.22%  "

# Row 7
$ws.Range("D7").Value = "1.008"
$ws.Range("E7").Value = "  +0.47%  "

# Row 8
$ws.Range("D8").Value = "0.06459"
$ws.Range("E8").Value = "  -1.78%  "

# Row 9
$ws.Range("D9").Value = "0.2561"
$ws.Range("E9").Value = "  -3.37%  "

# Row 10
$ws.Range("E10").Value = "  -3.77%  "

# Row 11
$ws.Range("D11").Value = "0.07646"
$ws.Range("E11").Value = "  +0.05%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.343"
$ws.Range("E12").Value = "  -5.15%  "

# Row 13
$ws.Range("D13").Value = "1.673.58"
$ws.Range("E13").Value = "  -1.96%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.900.54"
$ws.Range("E14").Value = "  -1.97%  "

# Row 15
$ws.Range("D15").Value = "0.5570"
$ws.Range("E15").Value = "  -2.63%  "

# Row 16
$ws.Range("D16").Value = "0.0₅8002"
$ws.Range("E16").Value = "  -2.04%  "

# Row 17
$ws.Range("D17").Value = "65.06"
$ws.Range("E17").Value = "  -3.47%  "

# Row 18
$ws.Range("D18").Value = "26.549.76"
$ws.Range("E18").Value = "  -2.46%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.52%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "210.41"
$ws.Range("E20").Value = "  -2.54%  "

# Row 21
$ws.Range("D21").Value = "4.452"
$ws.Range("E21").Value = "  -4.52%  "

# Row 22
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  -3.50%  "

# Row 23
$ws.Range("D23").Value = "5.885"
$ws.Range("E23").Value = "  -1.38%  "

# Row 24
$ws.Range("D24").Value = "1.008"
$ws.Range("E24").Value = "  +0.48%  "

# Row 25
$ws.Range("D25").Value = "142.63"
$ws.Range("E25").Value = "  +0.41%  "

# Row 26
$ws.Range("D26").Value = "1.724"
$ws.Range("E26").Value = "  -0.97%  "

# Row 27
$ws.Range("D27").Value = "0.1165"
$ws.Range("E27").Value = "  -4.35%  "

# Row 28
$ws.Range("D28").Value = "6.990"
$ws.Range("E28").Value = "  -3.49%  "

# Row 29
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  -3.73%  "

# Row 30
$ws.Range("D30").Value = "0.05225"
$ws.Range("E30").Value = "  -2.68%  "

# Row 31
$ws.Range("D31").Value = "1.265"
$ws.Range("E31").Value = "  -1.82%  "

# Row 32
$ws.Range("D32").Value = "3.348"
$ws.Range("E32").Value = "  -4.49%  "

# Row 33
$ws.Range("D33").Value = "3.202"
$ws.Range("E33").Value = "  -6.19%  "

# Row 34
$ws.Range("D34").Value = "1.581"
$ws.Range("E34").Value = "  -3.12%  "

# Row 35
$ws.Range("D35").Value = "2.758"
$ws.Range("E35").Value = "  -4.01%  "

# Row 36
$ws.Range("D36").Value = "2.386"
$ws.Range("E36").Value = "  -1.56%  "

# Row 37
$ws.Range("D37").Value = "0.9232"
$ws.Range("E37").Value = "  -2.34%  "

# Row 38
$ws.Range("D38").Value = "0.5683"
$ws.Range("E38").Value = "  -2.86%  "

# Row 39
$ws.Range("D39").Value = "1.161.35"
$ws.Range("E39").Value = "  +11.89%  "

# Row 40
$ws.Range("D40").Value = "0.01590"
$ws.Range("E40").Value = "  -2.44%  "

# Row 41
$ws.Range("D41").Value = "1.008"
$ws.Range("E41").Value = "  +0.47%  "

# Row 42
$ws.Range("D42").Value = "5.669"
$ws.Range("E42").Value = "  -3.03%  "

# Row 43
$ws.Range("D43").Value = "0.8284"
$ws.Range("E43").Value = "  -1.19%  "

# Row 44
$ws.Range("D44").Value = "100.01"
$ws.Range("E44").Value = "  -0.94%  "

# Row 45
$ws.Range("D45").Value = "1.808.73"
$ws.Range("E45").Value = "  -2.04%  "

# Row 46
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -4.78%  "

# Row 47
$ws.Range("D47").Value = "0.4490"
$ws.Range("E47").Value = "  +0.08%  "

# Row 48
$ws.Range("D48").Value = "55.67"
$ws.Range("E48").Value = "  -3.97%  "

# Row 49
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.15%  "

# Row 50
$ws.Range("D50").Value = "7.947"
$ws.Range("E50").Value = "  -1.40%  "

# Row 51
$ws.Range("D51").Value = "0.05146"
$ws.Range("E51").Value = "  -1.72%  "

# Restore column D to the workbook default style (drops the temporary text
# number format so the saved XML matches the original unstyled cells).
$priceRange.Style = "Normal"

